$d = $word.ActiveDocument

# --- 1. Heading: "Bananasplit Rev 7: 3.5mm" -> "Banana split " + bookmark + "Rev 7: 3.5mm"
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$r1 = $d.Range(0, 12)
$r1.Text = "Banana split "
$d.Bookmarks.Add("_GoBack", $d.Range(13, 13))

# --- 2. "This note explains..." paragraph: "Banansplit" -> "Banana Split"
$d.Content.Find.Execute("Banansplit", $true, $false, $false, $false, $false, $true, 1, $false, "Banana Split", 2)

# --- 3. "Earlier versions of the Bananasplit ha" -> "Earlier versions of the Banana Split ha"
$d.Content.Find.Execute("versions of the Bananasplit ha", $true, $false, $false, $false, $false, $true, 1, $false, "versions of the Banana Split ha", 2)

# --- 4. Table cell: "Original Bananasplit (and Novation, Arturia etc)" -> "Original Banana split (and Novation, Arturia etc)"
$d.Content.Find.Execute("Original Bananasplit", $true, $false, $false, $false, $false, $true, 1, $false, "Original Banana split", 2)

# --- 5. Table cell (bold): "And Bananasplit v7" -> "And Banana split v7"
$d.Content.Find.Execute("And Bananasplit v7", $true, $false, $false, $false, $false, $true, 1, $false, "And Banana split v7", 2)

# --- 6. Bold paragraph: "We have now revised the Bananasplit PCB to use" -> "We have now revised the Banana split PCB to use"
$d.Content.Find.Execute("We have now revised the Bananasplit PCB to use", $true, $false, $false, $false, $false, $true, 1, $false, "We have now revised the Banana split PCB to use", 2)

# --- 7. "do not unscew the bolts" -> "do not unscrew the bolts"
$d.Content.Find.Execute("unscew", $true, $false, $false, $false, $false, $true, 1, $false, "unscrew", 2)
